$d = $word.ActiveDocument

# 1. Change the date
$d.Content.Find.Execute("2024-03-06", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-07", 2) | Out-Null

# 2. Update the Introduction paragraph wording / citation
$d.Content.Find.Execute("Introduction goes here.", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction goes here, following", 2) | Out-Null
$d.Content.Find.Execute("Marrero et al. (2019)", $true, $false, $false, $false, $false, $true, 1, $false, "(Stockall et al. 2019; Neophytou et al. 2018)", 2) | Out-Null

# 3. Add a new "Consequently." paragraph (BodyText style) right after the Introduction paragraph
$introPara = $d.Paragraphs(10)
$introPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(11)
$newPara.Range.Text = "Consequently."
$newPara.Style = "BodyText"

# 4. Replace the Marrero et al. 2019 bibliography entry with the Neophytou et al.
#    2018 entry, and append a new Stockall et al. 2019 entry after it.
$refPara = $d.Paragraphs(43)
$startPos = $refPara.Range.Start
$endPos = $refPara.Range.End

# Wipe the paragraph's existing (mixed-formatting) content, keep the paragraph mark.
$wipe = $d.Range($startPos, $endPos - 1)
$wipe.Text = ""

$cur = $d.Range($startPos, $startPos)
$cur.InsertAfter("Neophytou, K., C. Manouilidou, L. Stockall, and A. Marantz. 2018.")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter(" ")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter("“Syntactic and Semantic Restrictions on Morphological Recomposition: MEG Evidence from Greek.”")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter(" ")
$cur = $d.Range($cur.End, $cur.End)
$italicStart = $cur.End
$cur.InsertAfter("Brain and Language")
$italicRange = $d.Range($italicStart, $cur.End)
$italicRange.Font.Italic = $true
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter(" ")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter("183 (August): 11–20.")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter(" ")
$cur = $d.Range($cur.End, $cur.End)
$linkStart = $cur.End
$cur.InsertAfter("https://doi.org/10.1016/j.bandl.2018.05.003")
$linkRange = $d.Range($linkStart, $cur.End)
$d.Hyperlinks.Add($linkRange, "https://doi.org/10.1016/j.bandl.2018.05.003") | Out-Null
$cur = $d.Range($refPara.Range.End - 1, $refPara.Range.End - 1)
$cur.InsertAfter(".")

# Add the new Stockall et al. 2019 bibliography entry as its own paragraph.
$refPara.Range.InsertParagraphAfter()
$stockPara = $d.Paragraphs(44)
$stockPara.Style = "Bibliography"

$sStart = $stockPara.Range.Start
$cur = $d.Range($sStart, $sStart)
$cur.InsertAfter("Stockall, Linnaea, Christina Manouilidou, Laura Gwilliams, Kyriaki Neophytou, and Alec Marantz. 2019.")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter(" ")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter("“Prefix Stripping Re-Re-Revisited: MEG Investigations of Morphological Decomposition and Recomposition.”")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter(" ")
$cur = $d.Range($cur.End, $cur.End)
$italicStart2 = $cur.End
$cur.InsertAfter("Frontiers in Psychology")
$italicRange2 = $d.Range($italicStart2, $cur.End)
$italicRange2.Font.Italic = $true
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter(" ")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter("10 (September).")
$cur = $d.Range($cur.End, $cur.End)
$cur.InsertAfter(" ")
$cur = $d.Range($cur.End, $cur.End)
$linkStart2 = $cur.End
$cur.InsertAfter("https://doi.org/10.3389/fpsyg.2019.01964")
$linkRange2 = $d.Range($linkStart2, $cur.End)
$d.Hyperlinks.Add($linkRange2, "https://doi.org/10.3389/fpsyg.2019.01964") | Out-Null
$cur = $d.Range($stockPara.Range.End - 1, $stockPara.Range.End - 1)
$cur.InsertAfter(".")
